# feat: add 2022-Q3 data
#
# Target layout (sheet order): 总计, 2022-Q3, 2022-Q2
#   - "总计" gains a 3rd row ("2022-Q2" copied down) and its existing row 2
#     is relabelled "2022-Q3".
#   - The existing "2022-Q2" sheet (rId2/sheetId2) is duplicated so the
#     duplicate (new rId/sheetId) preserves the untouched old Q2 fund data
#     and keeps the "2022-Q2" name.
#   - The original sheet (still rId2/sheetId2) is renamed "2022-Q3" and its
#     contents are replaced with the new quarter's fund data.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# --- 1. Duplicate the current "2022-Q2" sheet right after itself. The copy
#        keeps the old data untouched and will remain "2022-Q2". ---------
$q2Sheet.Copy([System.Reflection.Missing]::Value, $q2Sheet)

$newQ3Sheet = $wb.Worksheets.Item(2)
$oldQ2Copy = $wb.Worksheets.Item(3)

# Rename the still-original sheet first (it currently holds the name
# "2022-Q2", so the copy must wait its turn before reclaiming that name).
$newQ3Sheet.Name = "2022-Q3"
$oldQ2Copy.Name = "2022-Q2"

# --- 2. Replace the (now renamed) "2022-Q3" sheet's data with the new
#        quarter numbers, matching the "总计" sheet's header/cell style. ---
$newQ3Sheet.Range("B1").Value = "基金代码"
$newQ3Sheet.Range("C1").Value = "基金名称"
$newQ3Sheet.Range("D1").Value = "基金规模"
$newQ3Sheet.Range("E1").Value = "股票总仓位"
$newQ3Sheet.Range("F1").Value = "仓位占比"
$newQ3Sheet.Range("G1").Value = "持有市值(亿元)"
$newQ3Sheet.Range("H1").Value = "仓位排名"

$newQ3Sheet.Range("A2").Value = 0
$newQ3Sheet.Range("C2").Value = "工银全球精选股票（QDII）"
$newQ3Sheet.Range("H2").Value = 6

# These source values look numeric ("486002", "3.72", ...) but the original
# workbook stores them as plain text, not numbers. Force text via a
# scratch cell formatted "@" (Text), copy/paste its value+format into the
# target cells, then drop the scratch column so it leaves no trace in the
# used range.
$scratch = $newQ3Sheet.Range("Z1")
$scratch.NumberFormat = "@"
$xlPasteValuesAndNumberFormats = -4163

$scratch.Value = "486002"
$scratch.Copy()
$newQ3Sheet.Range("B2").PasteSpecial($xlPasteValuesAndNumberFormats)

$scratch.Value = "3.72"
$scratch.Copy()
$newQ3Sheet.Range("D2").PasteSpecial($xlPasteValuesAndNumberFormats)

$scratch.Value = "93.69"
$scratch.Copy()
$newQ3Sheet.Range("E2").PasteSpecial($xlPasteValuesAndNumberFormats)

$scratch.Value = "1.87"
$scratch.Copy()
$newQ3Sheet.Range("F2").PasteSpecial($xlPasteValuesAndNumberFormats)

$scratch.Value = "0.0696"
$scratch.Copy()
$newQ3Sheet.Range("G2").PasteSpecial($xlPasteValuesAndNumberFormats)

$scratch.EntireColumn.Delete()

# Match styling (bold/center/thin-border style used by "总计") and page
# margins on the rebuilt sheet.
$totalSheet.Range("B1:D1").Copy()
$newQ3Sheet.Range("B1:H1").PasteSpecial($xlPasteFormats)
$totalSheet.Range("A2").Copy()
$newQ3Sheet.Range("A2").PasteSpecial($xlPasteFormats)

$newQ3Sheet.PageSetup.LeftMargin = 54
$newQ3Sheet.PageSetup.RightMargin = 54
$newQ3Sheet.PageSetup.TopMargin = 72
$newQ3Sheet.PageSetup.BottomMargin = 72
$newQ3Sheet.PageSetup.HeaderMargin = 36
$newQ3Sheet.PageSetup.FooterMargin = 36

# --- 3. Update "总计": row 2 now reports "2022-Q3", and a new row 3 is
#        added for the (still present) "2022-Q2" figures. -----------------
$totalSheet.Range("B2").Value = "2022-Q3"

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.07000000000000001

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial($xlPasteFormats)
